# stats de structure #27 : Ajout d'un onglet "à la date de l'extraction"
#
# The "Chiffres" sheet (a per-structure stats table driven by a templating
# engine via shared-string placeholders) is duplicated so there is one
# version that lists effectifs as of the end of the scouting season, and a
# new version that lists them as of the report's generation date.

$wb = $excel.ActiveWorkbook

$wsChiffres = $wb.Worksheets.Item("Chiffres")

# Duplicate "Chiffres" -> the copy is inserted right after the original and
# becomes the active sheet.
$wsChiffres.Copy([System.Reflection.Missing]::Value, $wsChiffres)
$wsParDate = $wb.Worksheets.Item(2)

# Rename both tabs to reflect their respective data sources.
$wsChiffres.Name = "Par saison (fin de saison)"
$wsParDate.Name = "Par saison (date de génération)"

# Point the "fin de saison" sheet's templating placeholder at the
# end-of-season effectifs collection (the copy keeps the original
# "${effectifs}" placeholder).
$wsChiffres.Range("A5").Formula = '<jt:forEach items="${effectifs_findannee}" var="effectif">${effectif.groupe}'

# Leave the original sheet as the selected/active tab with A5 selected.
$null = $wsChiffres.Activate()
$null = $wsChiffres.Range("A5").Select()
